# Applies the "Dudas" row-8 addition described by the commit:
#   - a new row (#6) is added to the Q&A table with a new question about
#     list-printing formatting
#   - the active selection moves from F18 to B9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Build row 8 as a copy of row 7's formatting (border/font/fill for
#        the #, Pregunta and Solucionada columns), then set its values.
$null = $ws.Range("A7:C7").Copy()
$null = $ws.Range("A8:C8").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Se puede cambiar el formato del listado para que tenga un espacio antes y despues`ndel "" | "" y tmb tuvimos que inicializar el valorString ya que sino imprimia null al`ninicio del listado"

# Row height follows the same auto-sized pattern as the other question
# rows (~19.5 per wrapped line once the text is wrapped at column B's
# width) -- for this 3-line question that comes out to 57pt.
$ws.Rows.Item(8).RowHeight = 57

# --- 2. Update the saved selection to match the authored file (B9).
$null = $ws.Range("B9").Select()
